$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The duplicated/stray "preview" values that had been pasted into column L
# (rows 14-19) are removed, leaving the cells' formatting untouched.
$ws.Range("L14:L19").ClearContents()

# Restore the saved selection to L13.
$ws.Range("L13").Select()
